$d = $word.ActiveDocument
$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) "DestroyABarrier" paragraph: drop the yellow highlight that was applied
#    both to the paragraph mark (pPr/rPr) and to the run of text.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("DestroyABarrier") | Out-Null
$para1 = $rng1.Paragraphs(1)
$xml1 = '<w:p xmlns:w="' + $wordNs + '"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>DestroyABarrier</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$para1.Range.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new bullet "la ficha entra en el pasillo" (level 4 of the same
#    list) right before the "EnterFinalSquare" paragraph.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("EnterFinalSquare") | Out-Null
$para2 = $rng2.Paragraphs(1)
$para2.Range.InsertParagraphBefore()
$para2.Range.Text = "la ficha entra en el pasillo"
$para2.Range.ListFormat.ListLevelNumber = 5

# ---------------------------------------------------------------------------
# 3) Rework the trailing notes block: the old "Cambiar funcion expand..."
#    paragraph (which carried the single-underline paragraph mark) is split
#    into four paragraphs - the original text now plain, followed by three
#    new notes, with the underline paragraph-mark formatting moved to the
#    2nd and 4th of them.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Cambiar funcion expand") | Out-Null
$para3 = $rng3.Paragraphs(1)
$xml3 = (
  '<w:p xmlns:w="' + $wordNs + '">' +
    '<w:r><w:t xml:space="preserve">Cambiar </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>funcion</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>expand</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> para que pase menos </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>parametros</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> &#191;?</w:t></w:r>' +
  '</w:p>' +
  '<w:p xmlns:w="' + $wordNs + '">' +
    '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">la </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>information</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> del </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>state</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> es la partida</w:t></w:r>' +
  '</w:p>' +
  '<w:p xmlns:w="' + $wordNs + '">' +
    '<w:r><w:t xml:space="preserve">en el </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>effect</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> hay que poner el rating y el movimiento de las fichas!!!!</w:t></w:r>' +
  '</w:p>' +
  '<w:p xmlns:w="' + $wordNs + '">' +
    '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">104 significa que </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> en casa</w:t></w:r>' +
  '</w:p>'
)
$para3.Range.InsertXML($xml3) | Out-Null

Write-Host "edits applied"
